$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for all existing data rows (2-223)
# from 45181 (2023-09-12) to 45182 (2023-09-13).
$ws.Range("C2:C223").Value = 45182

# Row 223 gains an explicit row height (ht="15" customHeight="1").
$ws.Rows.Item(223).RowHeight = 15

# Append the new record as row 224.
$ws.Range("A224").Value = "A 42304-2023"
$ws.Range("B224").NumberFormat = "YYYY-MM-DD"
$ws.Range("B224").Value = 45180
$ws.Range("C224").NumberFormat = "YYYY-MM-DD"
$ws.Range("C224").Value = 45182
$ws.Range("D224").Value = "ÖREBRO LÄN"
$ws.Range("E224").Value = "LJUSNARSBERG"
$ws.Range("G224").Value = 27
$ws.Range("H224").Value = 0
$ws.Range("I224").Value = 0
$ws.Range("J224").Value = 0
$ws.Range("K224").Value = 0
$ws.Range("L224").Value = 0
$ws.Range("M224").Value = 0
$ws.Range("N224").Value = 0
$ws.Range("O224").Value = 0
$ws.Range("P224").Value = 0
$ws.Range("Q224").Value = 0
$ws.Range("R224").WrapText = $true
$ws.Range("R224").Value = ""
